$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Avatar', ['Token Creature — Avatar', 'This creature’s power and toughness are each equal to your life total.', '*/*'])"
$ws.Range("A3").Value = "('Beast', ['Token Creature — Beast', '3/3'])"
$ws.Range("A4").Value = "('Bird', ['Token Creature — Bird', 'Flying', '3/3'])"
$ws.Range("A5").Value = "('Ooze', ['Token Creature — Ooze', 'When this creature dies, create two 1/1 green Ooze creature tokens.', '2/2'])"
$ws.Range("A6").Value = "('Zombie', ['Token Creature — Zombie', '2/2'])"

$ws.Rows("7:22").Delete()
